# Update "想去人数" (interested-attendee count) figures in column F across
# the sheets that list individual event rows. Sheet "本地生活" has no
# changed rows, so it is left untouched.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2  = 176
        5  = 946
        6  = 5032
        7  = 401
        8  = 578
        9  = 877
        16 = 1649
        18 = 733
        21 = 267
        23 = 120
        24 = 1041
        27 = 2001
        28 = 155
        31 = 212
        36 = 260
        37 = 572
        39 = 29
    }
    "演出" = @{
        6 = 97
    }
    "全部类型" = @{
        3  = 176
        6  = 946
        8  = 5032
        9  = 401
        10 = 578
        13 = 877
        16 = 97
        23 = 1649
        25 = 733
        28 = 267
        31 = 120
        32 = 1041
        34 = 2001
        35 = 155
        38 = 212
        42 = 260
        43 = 572
        45 = 29
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
